$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update photograph path cells (D2:D4) -> single common path
$ws.Range("D2").Value = "C:\Users\Jama\Pictures\testng-tutorial.png"
$ws.Range("D3").Value = "C:\Users\Jama\Pictures\testng-tutorial.png"
$ws.Range("D4").Value = "C:\Users\Jama\Pictures\testng-tutorial.png"

# Update username cells (E2:E4)
$ws.Range("E2").Value = "anasule001234567"
$ws.Range("E3").Value = "blakenailya001234567"
$ws.Range("E4").Value = "mikeaj001234567"

# Resize columns D and E to fit the new content
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws.Columns.Item(5).ColumnWidth = 18.6

# Update selected cell in the view
$ws.Range("F10").Select()
